$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the "land/mortgage liability" record. The "★" marker in B2
# is dropped entirely (cell becomes blank, keeping its border/style).
$ws.Range("B2").ClearContents()

# Tidy up stray internal spaces in the creditor-address, balance and date
# strings (C2/G2 - 蘇震清/貸款 - are left untouched).
$ws.Range("D2").Value = "臺灣新光商業銀行五常分行臺北市中山區龍江路356巷76號"

# Balance loses its thousands separators; keep it as text (not a number)
# by forcing Excel's text quote-prefix, then drop back to the plain
# "Normal" style so no stray numeric formatting style is introduced.
$ws.Range("E2").Value = "'1468576"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "100年05月16日"
